$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects("Tabela9")
$lo.Resize($ws.Range("I3:L23")) | Out-Null
$ws.Range("I22").Value = "executáve(jar)"
$ws.Range("I22").HorizontalAlignment = -4131 # xlLeft
$ws.Range("I22").Borders.Item(10).Weight = -4138 # right edge, medium
